# Adds new "bunnies" tournament definitions (time attack / time limit / normal)
# to the "tournaments" sheet and a new "kill bunnies" global-event quest
# definition to the "quests" sheet, mirroring the formatting of existing rows.

$wb = $excel.ActiveWorkbook
$wsTournaments = $wb.Worksheets.Item("tournaments")
$wsQuests = $wb.Worksheets.Item("quests")

# ---------------------------------------------------------------------------
# 1) "tournaments" sheet: append rows 170-172, copying the formatting of an
#    existing, same-shaped row (row 21: A=9 B=7 C=9 D=9 E=7 F=9 G=11) and then
#    clearing column B's style back to the workbook default, matching the
#    un-styled <c r="B.." t="s"> cells used for these new rows.
# ---------------------------------------------------------------------------
$wsTournaments.Range("A21:G21").Copy($wsTournaments.Range("A170:G170"))
$wsTournaments.Range("A21:G21").Copy($wsTournaments.Range("A171:G171"))
$wsTournaments.Range("A21:G21").Copy($wsTournaments.Range("A172:G172"))

$wsTournaments.Range("B170").Style = "Normal"
$wsTournaments.Range("B171").Style = "Normal"
$wsTournaments.Range("B172").Style = "Normal"

# Row 172 first so "...TIME_ATTACK_BUNNIES" becomes the first newly added
# shared string, then row 171, then row 170 -- this reproduces the exact
# shared-string insertion order of the authored workbook.
$wsTournaments.Range("B172").Value = "TID_EVENT_TOURNAMENT_KILL_TIME_ATTACK_BUNNIES"
$wsTournaments.Range("C172").Value = "kill"
$wsTournaments.Range("D172").Value = 1

$wsTournaments.Range("B171").Value = "TID_EVENT_TOURNAMENT_KILL_TIME_LIMIT_BUNNIES"
$wsTournaments.Range("C171").Value = "kill"
$wsTournaments.Range("D171").Value = 2

$wsTournaments.Range("B170").Value = "TID_EVENT_TOURNAMENT_KILL_NORMAL_BUNNIES"
$wsTournaments.Range("C170").Value = "kill"
$wsTournaments.Range("D170").Value = 0

$wsTournaments.Range("E170").Value = "FlyingBunny;easter_bunny"
$wsTournaments.Range("E171").Value = "FlyingBunny;easter_bunny"
$wsTournaments.Range("E172").Value = "FlyingBunny;easter_bunny"

$wsTournaments.Range("G170").Value = "icon_bunnies"
$wsTournaments.Range("G171").Value = "icon_bunnies"
$wsTournaments.Range("G172").Value = "icon_bunnies"

$wsTournaments.Range("F170").ClearContents()
$wsTournaments.Range("F171").ClearContents()
$wsTournaments.Range("F172").ClearContents()

$wsTournaments.Range("A170").Value = "<Definition>"
$wsTournaments.Range("A171").Value = "<Definition>"
$wsTournaments.Range("A172").Value = "<Definition>"

# ---------------------------------------------------------------------------
# 2) "quests" sheet: append row 58, copying the formatting of an existing,
#    same-shaped row (row 6: A=9 B=7 C=9 D=7 E=9 F=12).
# ---------------------------------------------------------------------------
$wsQuests.Range("A6:F6").Copy($wsQuests.Range("A58:F58"))

$wsQuests.Range("B58").Value = "TID_GLOBAL_EVENT_KILL_BUNNIES"
$wsQuests.Range("C58").Value = "kill"
$wsQuests.Range("D58").Value = "FlyingBunny;easter_bunny"
$wsQuests.Range("F58").Value = "icon_bunnies"
$wsQuests.Range("E58").ClearContents()
$wsQuests.Range("A58").Value = "<Definition>"

# ---------------------------------------------------------------------------
# 3) Restore the view/selection state recorded in the authored workbook.
# ---------------------------------------------------------------------------
$wsTournaments.Activate()
[void]$wsTournaments.Range("E172").Select()
try { $excel.ActiveWindow.ScrollRow = 145 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 2 } catch {}

$wsQuests.Activate()
[void]$wsQuests.Range("D59").Select()
try { $excel.ActiveWindow.ScrollRow = 37 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
